$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E5").Value = 9
$ws.Range("G5").Value = 2250

$ws.Range("E6").Value = 90
$ws.Range("G6").Value = 2250
